$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 5) for "2021年", replicating the formatting of the
# preceding data row (row 4) so the new row label cell (A5) keeps the
# same bold/centered/bordered style used by the other year-label cells.
$ws.Range("A4:M4").Copy($ws.Range("A5:M5"))

$ws.Range("A5").Value = "2021年"
$ws.Range("B5").Value = 43.396
$ws.Range("C5").Value = 35.581
$ws.Range("D5").Value = 18.997
$ws.Range("E5").Value = 21.507
$ws.Range("F5").Value = 37.057
$ws.Range("G5").Value = 33.813
$ws.Range("H5").Value = 37.97
$ws.Range("I5").Value = 35.425
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 37.017
$ws.Range("L5").Value = 63.463
$ws.Range("M5").Value = 11.534
